$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 926.1892
$ws.Range("I15").Value = 926.1892
$ws.Range("K15").Value = 2778.5676
$ws.Range("M15").Value = -2609.5676
$ws.Range("H107").Value = 36622.723
$ws.Range("I107").Value = 45842.74
$ws.Range("K107").Value = 45842.74
$ws.Range("M107").Value = -43922.74
$ws.Range("H137").Value = 3692.087
$ws.Range("I137").Value = 1979.091
$ws.Range("J137").Value = 5262.3335
$ws.Range("K137").Value = 5937.272999999999
$ws.Range("L137").Value = 15787.0005
$ws.Range("M137").Value = -3387.272999999999
$ws.Range("N137").Value = -20887.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2365.5283
$ws.Range("I32").Value = 2365.5283
$ws.Range("K32").Value = 2365.5283
$ws.Range("M32").Value = -2078.5283
$ws.Range("H122").Value = 4945.154
$ws.Range("I122").Value = 4392.647
$ws.Range("K122").Value = 13177.941
$ws.Range("M122").Value = -10727.941

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H134").Value = 70872.125
$ws.Range("I134").Value = 8905.4
$ws.Range("J134").Value = 174150
$ws.Range("K134").Value = 26716.2
$ws.Range("L134").Value = 522450
$ws.Range("M134").Value = -24181.2
$ws.Range("N134").Value = -527520

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3753
$ws.Range("I16").Value = 2999
$ws.Range("K16").Value = 2999
$ws.Range("M16").Value = -2712
$ws.Range("H31").Value = 3906.7334
$ws.Range("I31").Value = 2121.5557
$ws.Range("K31").Value = 2121.5557
$ws.Range("M31").Value = -1826.5557
$ws.Range("H34").Value = 3906.7334
$ws.Range("I34").Value = 2121.5557
$ws.Range("K34").Value = 2121.5557
$ws.Range("M34").Value = -1919.5557
$ws.Range("H58").Value = 347743.56
$ws.Range("I58").Value = 528380.9399999999
$ws.Range("J58").Value = 4532.5
$ws.Range("K58").Value = 528380.9399999999
$ws.Range("L58").Value = 4532.5
$ws.Range("M58").Value = -528177.9399999999
$ws.Range("N58").Value = -4938.5
$ws.Range("H99").Value = 6101.909
$ws.Range("I99").Value = 4352.25
$ws.Range("K99").Value = 4352.25
$ws.Range("M99").Value = -2854.25
$ws.Range("H105").Value = 1957.3636
$ws.Range("I105").Value = 1838.1
$ws.Range("K105").Value = 1838.1
$ws.Range("M105").Value = -91.09999999999991
$ws.Range("H113").Value = 3753
$ws.Range("I113").Value = 2999
$ws.Range("K113").Value = 2999
$ws.Range("M113").Value = -829
$ws.Range("H126").Value = 6101.909
$ws.Range("I126").Value = 4352.25
$ws.Range("K126").Value = 13056.75
$ws.Range("M126").Value = -10586.75
$ws.Range("H132").Value = 3303
$ws.Range("I132").Value = 3363.8
$ws.Range("K132").Value = 10091.4
$ws.Range("M132").Value = -7561.400000000001
$ws.Range("H134").Value = 530659.4
$ws.Range("J134").Value = 1433449.9
$ws.Range("L134").Value = 4300349.699999999
$ws.Range("N134").Value = -4305419.699999999
$ws.Range("H136").Value = 347743.56
$ws.Range("I136").Value = 528380.9399999999
$ws.Range("J136").Value = 4532.5
$ws.Range("K136").Value = 1585142.82
$ws.Range("L136").Value = 13597.5
$ws.Range("M136").Value = -1582592.82
$ws.Range("N136").Value = -18697.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2271137.8
$ws.Range("J4").Value = 12210030
$ws.Range("L4").Value = 36630090
$ws.Range("N4").Value = -36630314
$ws.Range("H113").Value = 7408209
$ws.Range("I113").Value = 37037036
$ws.Range("J113").Value = 1002.5
$ws.Range("K113").Value = 111111108
$ws.Range("L113").Value = 3007.5
$ws.Range("M113").Value = -111108938
$ws.Range("N113").Value = -7347.5
$ws.Range("H124").Value = 18000
$ws.Range("I124").Value = 18000
$ws.Range("K124").Value = 54000
$ws.Range("M124").Value = -49090
$ws.Range("H130").Value = 2499.6667
$ws.Range("I130").Value = 2499.5
$ws.Range("J130").Value = 2500
$ws.Range("K130").Value = 7498.5
$ws.Range("L130").Value = 7500
$ws.Range("M130").Value = -2478.5
$ws.Range("N130").Value = -17540
$ws.Range("H131").Value = 3665.4783
$ws.Range("I131").Value = 1270.2307
$ws.Range("J131").Value = 6779.3
$ws.Range("K131").Value = 3810.6921
$ws.Range("L131").Value = 20337.9
$ws.Range("M131").Value = 1229.3079
$ws.Range("N131").Value = -30417.9
$ws.Range("H132").Value = 3099.0908
$ws.Range("I132").Value = 2013.5714
$ws.Range("J132").Value = 4998.75
$ws.Range("K132").Value = 18122.1426
$ws.Range("L132").Value = 44988.75
$ws.Range("M132").Value = -15592.1426
$ws.Range("N132").Value = -50048.75
$ws.Range("H134").Value = 3311.5715
$ws.Range("I134").Value = 2864.4546
$ws.Range("J134").Value = 4951
$ws.Range("K134").Value = 8593.363799999999
$ws.Range("L134").Value = 14853
$ws.Range("M134").Value = -3523.363799999999
$ws.Range("N134").Value = -24993
$ws.Range("H140").Value = 3360.4666
$ws.Range("I140").Value = 3360.4666
$ws.Range("K140").Value = 10081.3998
$ws.Range("M140").Value = -4901.399800000001
$ws.Range("H141").Value = 2850.8
$ws.Range("I141").Value = 2850.8
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 8552.400000000001
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -3372.400000000001
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 71435870
$ws.Range("I70").Value = 8047.875
$ws.Range("K70").Value = 8047.875
$ws.Range("M70").Value = -7777.875
$ws.Range("H73").Value = 71435870
$ws.Range("I73").Value = 8047.875
$ws.Range("K73").Value = 8047.875
$ws.Range("M73").Value = -7111.875
$ws.Range("H102").Value = 3201.3462
$ws.Range("I102").Value = 2644.8333
$ws.Range("K102").Value = 2644.8333
$ws.Range("M102").Value = -1022.8333
$ws.Range("H113").Value = 376766.34
$ws.Range("I113").Value = 533699.4399999999
$ws.Range("K113").Value = 533699.4399999999
$ws.Range("M113").Value = -531529.4399999999
$ws.Range("H122").Value = 926843.5600000001
$ws.Range("I122").Value = 1008647.56
$ws.Range("K122").Value = 3025942.68
$ws.Range("M122").Value = -3023492.68
$ws.Range("H132").Value = 231746.8
$ws.Range("I132").Value = 912114.9399999999
$ws.Range("J132").Value = 29475.19
$ws.Range("K132").Value = 2736344.82
$ws.Range("L132").Value = 88425.56999999999
$ws.Range("M132").Value = -2733814.82
$ws.Range("N132").Value = -93485.56999999999
$ws.Range("H133").Value = 69983.336
$ws.Range("J133").Value = 69983.336
$ws.Range("L133").Value = 69983.336
$ws.Range("N133").Value = -80103.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 737.1818
$ws.Range("I22").Value = 836.25
$ws.Range("J22").Value = 473
$ws.Range("K22").Value = 836.25
$ws.Range("L22").Value = 473
$ws.Range("M22").Value = -541.25
$ws.Range("N22").Value = -1063
$ws.Range("H27").Value = 737.1818
$ws.Range("I27").Value = 836.25
$ws.Range("J27").Value = 473
$ws.Range("K27").Value = 836.25
$ws.Range("L27").Value = 473
$ws.Range("M27").Value = -729.25
$ws.Range("N27").Value = -687
$ws.Range("H43").Value = 1000000
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H46").Value = 4861.737
$ws.Range("J46").Value = 5944.3335
$ws.Range("L46").Value = 5944.3335
$ws.Range("N46").Value = -6320.3335
$ws.Range("H61").Value = 7355.75
$ws.Range("I61").Value = 5269
$ws.Range("K61").Value = 5269
$ws.Range("M61").Value = -5067
$ws.Range("H100").Value = 114439.9
$ws.Range("I100").Value = 264499.75
$ws.Range("J100").Value = 14400
$ws.Range("K100").Value = 264499.75
$ws.Range("L100").Value = 14400
$ws.Range("M100").Value = -263958.75
$ws.Range("N100").Value = -15482
$ws.Range("H113").Value = 7355.75
$ws.Range("I113").Value = 5269
$ws.Range("K113").Value = 5269
$ws.Range("M113").Value = -3099

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1120.2963
$ws.Range("I113").Value = 1014.0526
$ws.Range("K113").Value = 3042.1578
$ws.Range("M113").Value = -872.1578
$ws.Range("H132").Value = 32917.15
$ws.Range("I132").Value = 2551
$ws.Range("J132").Value = 117267.555
$ws.Range("K132").Value = 7653
$ws.Range("L132").Value = 351802.665
$ws.Range("M132").Value = -5123
$ws.Range("N132").Value = -356862.665
$ws.Range("H136").Value = 302881.62
$ws.Range("I136").Value = 359349.22
$ws.Range("J136").Value = 171123.92
$ws.Range("K136").Value = 1078047.66
$ws.Range("L136").Value = 513371.76
$ws.Range("M136").Value = -1075497.66
$ws.Range("N136").Value = -518471.76
